$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (rich text runs) ---
$ws.Range("A8").Characters(21, 2).Text = "12"
$ws.Range("C9").Characters(28, 9).Text = "3/18/2024"
$ws.Range("C9").Characters(49, 9).Text = "3/24/2024"

# --- Plain numeric cell updates ---
$ws.Range("M15").Value = -50
$ws.Range("N15").Value = -60
$ws.Range("C16").Value = 2
$ws.Range("E16").Value = 100
$ws.Range("F16").Value = 11
$ws.Range("G16").Value = 4
$ws.Range("H16").Value = 175
$ws.Range("I16").Value = 29
$ws.Range("J16").Value = 24
$ws.Range("K16").Value = 20.833333333333
$ws.Range("L16").Value = 26.086956521739
$ws.Range("M16").Value = -14.705882352941
$ws.Range("N16").Value = -84.15300546448
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 6
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 12
$ws.Range("G17").Value = 19
$ws.Range("H17").Value = -36.842105263157
$ws.Range("I17").Value = 36
$ws.Range("J17").Value = 45
$ws.Range("K17").Value = -20
$ws.Range("L17").Value = 16.129032258064
$ws.Range("M17").Value = 44
$ws.Range("N17").Value = -38.983050847457
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = 7
$ws.Range("H18").Value = -85.714285714285
$ws.Range("J18").Value = 23
$ws.Range("K18").Value = -43.478260869565
$ws.Range("N18").Value = -90.845070422535
$ws.Range("C19").Value = 4
$ws.Range("D19").Value = 5
$ws.Range("E19").Value = -20
$ws.Range("F19").Value = 22
$ws.Range("G19").Value = 29
$ws.Range("H19").Value = -24.137931034482
$ws.Range("I19").Value = 80
$ws.Range("J19").Value = 87
$ws.Range("K19").Value = -8.045977011494
$ws.Range("L19").Value = 3.896103896103
$ws.Range("M19").Value = 50.943396226415
$ws.Range("N19").Value = -48.051948051948
$ws.Range("D20").Value = 3
$ws.Range("G20").Value = 8
$ws.Range("H20").Value = -87.5
$ws.Range("J20").Value = 21
$ws.Range("K20").Value = -76.190476190476
$ws.Range("N20").Value = -95.238095238095
$ws.Range("C21").Value = 12
$ws.Range("D21").Value = 17
$ws.Range("E21").Value = -29.411764705882
$ws.Range("F21").Value = 47
$ws.Range("G21").Value = 67
$ws.Range("H21").Value = -29.850746268656
$ws.Range("I21").Value = 165
$ws.Range("J21").Value = 202
$ws.Range("K21").Value = -18.316831683168
$ws.Range("L21").Value = -13.612565445026
$ws.Range("M21").Value = 23.134328358209
$ws.Range("N21").Value = -74.693251533742
$ws.Range("G22").Value = 5
$ws.Range("H22").Value = -100
$ws.Range("D23").Value = 3
$ws.Range("E23").Value = -33.333333333333
$ws.Range("F23").Value = 10
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 35
$ws.Range("J23").Value = 35
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 29.629629629629
$ws.Range("M23").Value = 66.666666666666
$ws.Range("D24").Value = 9
$ws.Range("E24").Value = -22.222222222222
$ws.Range("F24").Value = 34
$ws.Range("G24").Value = 28
$ws.Range("H24").Value = 21.428571428571
$ws.Range("I24").Value = 87
$ws.Range("J24").Value = 104
$ws.Range("K24").Value = -16.346153846153
$ws.Range("L24").Value = -20.90909090909
$ws.Range("M24").Value = -7.446808510638
$ws.Range("D25").Value = 3
$ws.Range("E25").Value = -66.666666666666
$ws.Range("F25").Value = 3
$ws.Range("H25").Value = -57.142857142857
$ws.Range("I25").Value = 15
$ws.Range("J25").Value = 36
$ws.Range("K25").Value = -58.333333333333
$ws.Range("L25").Value = -57.142857142857
$ws.Range("C26").Value = 7
$ws.Range("D26").Value = 11
$ws.Range("E26").Value = -36.363636363636
$ws.Range("F26").Value = 24
$ws.Range("G26").Value = 24
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 69
$ws.Range("J26").Value = 57
$ws.Range("K26").Value = 21.052631578947
$ws.Range("L26").Value = 13.11475409836
$ws.Range("M26").Value = 25.454545454545
$ws.Range("L27").Value = -20
$ws.Range("F28").Value = 2
$ws.Range("G28").Value = 3
$ws.Range("H28").Value = -33.333333333333
$ws.Range("I28").Value = 4
$ws.Range("K28").Value = -66.666666666666
$ws.Range("L28").Value = -33.333333333333

# --- Cells changing from numeric to text placeholder "0" (style 14) ---
$ws.Range("D22").Value = "'0"
$ws.Range("D31").Value = "'0"
$ws.Range("F22").Value = "'0"
$ws.Range("G15").Value = "'0"
$ws.Range("G27").Value = "'0"

# --- Cells changing from numeric to text placeholder "***.*" (style 14) ---
$ws.Range("E22").Value = "'***.*"
$ws.Range("E31").Value = "'***.*"
$ws.Range("H15").Value = "'***.*"
$ws.Range("H27").Value = "'***.*"

# Fix style for newly-text placeholder cells to match style 14 (copy format from C14)
$ws.Range("C14").Copy()
$ws.Range("D22,D31,E22,E31,F22,G15,G27,H15,H27").PasteSpecial(-4122)

# --- Cells changing from text placeholder to numeric (style 16) ---
$ws.Range("C25").Value = 1
$ws.Range("C28").Value = 1
$ws.Range("C16").Copy()
$ws.Range("C25,C28").PasteSpecial(-4122)

$excel.CutCopyMode = 0